$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.74
$ws.Range("A10").Value = -21.552
$ws.Range("A12").Value = -21.621
$ws.Range("C12").Value = -12.72
$ws.Range("D12").Value = -8.138999999999999
$ws.Range("D13").Value = -7.792
$ws.Range("C17").Value = -13.36
$ws.Range("A18").Value = -22.004
$ws.Range("D21").Value = -7.931
$ws.Range("C26").Value = -13.033
$ws.Range("C27").Value = -13.315
$ws.Range("C28").Value = -12.692
$ws.Range("D36").Value = -7.346000000000001
$ws.Range("A37").Value = -19.841
$ws.Range("C37").Value = -12.467
$ws.Range("D38").Value = -7.674999999999999
$ws.Range("D41").Value = -8.275
$ws.Range("D52").Value = -8.084999999999999
$ws.Range("A55").Value = -21.859
$ws.Range("C65").Value = -12.1
$ws.Range("D67").Value = -7.709000000000001
$ws.Range("A68").Value = -21.536
$ws.Range("C73").Value = -12.685
$ws.Range("A77").Value = -20.651
$ws.Range("A78").Value = -20.15
$ws.Range("C84").Value = -13.113
$ws.Range("C85").Value = -12.082
$ws.Range("D89").Value = -8.183
$ws.Range("C93").Value = -11.417
$ws.Range("C95").Value = -12.014
$ws.Range("D95").Value = -7.580000000000001
$ws.Range("C98").Value = -13.12
$ws.Range("C99").Value = -11.299
$ws.Range("C101").Value = -12.436
$ws.Range("D105").Value = -7.834000000000001
